$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 151-157 appended below the existing data (which ends at row 150).
# Column A: Date-Hora (date/time serial, formatted like the existing date column)
# Column B: Chuva (mm) - numeric 0, right aligned like existing B column
# Column C: Nivel (cm) - plain numeric value, default (general) style

$newRows = @(
    @{ Row = 151; A = 45219.458912037036; B = 0; C = 4861.8 },
    @{ Row = 152; A = 45219.500578703701; B = 0; C = 4861.3 },
    @{ Row = 153; A = 45219.542245370372; B = 0; C = 4862.5 },
    @{ Row = 154; A = 45219.583912037036; B = 0; C = 4861.8999999999996 },
    @{ Row = 155; A = 45219.625578703701; B = 0; C = 4861.8 },
    @{ Row = 156; A = 45219.667245370372; B = 0; C = 4861.6000000000004 },
    @{ Row = 157; A = 45219.708912037036; B = 0; C = 4860 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A: date/time value, formatted the same way as the rest of column A
    $cellA = $ws.Range("A$rowNum")
    $cellA.Value = $r.A
    $cellA.NumberFormat = "m/d/yy h:mm"

    # Column B: numeric 0, right aligned (matches style used by B2:B150)
    $cellB = $ws.Range("B$rowNum")
    $cellB.Value = $r.B
    $cellB.HorizontalAlignment = -4152

    # Column C: plain numeric value with default/general style (no explicit formatting)
    $cellC = $ws.Range("C$rowNum")
    $cellC.Value = $r.C
    $cellC.Style = "Normal"
}

# Update the selected cell shown in the sheet view to F153
$ws.Range("F153").Select()
